$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 74, shifting existing rows 74..144 down to 75..145.
$ws.Rows(74).Insert()

# Populate the newly inserted row 74 with the new weekly price observation.
$ws.Cells.Item(74, 1).Value  = 3
$ws.Cells.Item(74, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(74, 3).Value  = "Coquimbo"
$ws.Cells.Item(74, 4).Value  = 44512
$ws.Cells.Item(74, 5).Value  = 5
$ws.Cells.Item(74, 6).Value  = 100112010
$ws.Cells.Item(74, 7).Value  = "Achicoria"
$ws.Cells.Item(74, 8).Value  = "Sin especificar"
$ws.Cells.Item(74, 9).Value  = "Primera"
$ws.Cells.Item(74, 10).Value = 113
$ws.Cells.Item(74, 11).Value = 6000
$ws.Cells.Item(74, 12).Value = 6500
$ws.Cells.Item(74, 13).Value = 6243
$ws.Cells.Item(74, 14).Value = "`$/caja 16 unidades"
$ws.Cells.Item(74, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(74, 16).Value = 390
$ws.Cells.Item(74, 17).Value = 16
$ws.Cells.Item(74, 18).Value = "Hortaliza"
